# Update "想去人数" (F) and "最低票价" (G) figures for the matching rows on
# both the "展览" and "全部类型" sheets (they carry the same data).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1081
    $ws.Range("G2").Value = "不可售"

    $ws.Range("F5").Value = 4638

    $ws.Range("F7").Value = 389

    $ws.Range("F9").Value = 916

    $ws.Range("F11").Value = 1090

    $ws.Range("F13").Value = 588

    $ws.Range("F15").Value = 20

    $ws.Range("F16").Value = 269
}
